$d = $word.ActiveDocument

# Locate the list-item paragraph whose sole text is "Density for histograms"
# and remove it entirely (including its paragraph mark), leaving the
# surrounding list items untouched.
foreach ($p in $d.Paragraphs) {
    $text = $p.Range.Text.Trim("`r`a`n")
    if ($text -eq "Density for histograms") {
        $p.Range.Delete()
        break
    }
}
